$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New connector row: chipdip link added in column D, row 13 (same row as the
# "SIM Molex 475530001 Экранированный" / molex entry), mirroring the style
# already used for the other hyperlink cells in column C.
$ws.Range("D13").Value = "chipdip"
$ws.Hyperlinks.Add($ws.Range("D13"), "https://www.chipdip.ru/product/txb0108pwr-2") | Out-Null
$ws.Range("D13").Style = $ws.Range("C13").Style

# New column D needs an explicit width like the existing B/C columns.
$ws.Columns.Item(4).ColumnWidth = 20.33

# Selection moved down a few rows.
$ws.Range("C20").Select() | Out-Null
